# Apply the "finish updating lesson 5 extra practice assignments" edit.
#
# The change touches raw OOXML structure (tblLook flag, removed header-row
# cnfStyle/tcPr overrides moved into a new tblStylePr firstRow conditional
# format, numbering abstractNum indent/tab cleanup, and a bCs sibling for an
# existing bold run) that doesn't map onto a single high-level Word object
# model call, so we edit the package XML directly via the WordOpenXML
# round-trip property and write it back.

$d = $word.ActiveDocument
$xml = $d.WordOpenXML

# 1. Add <w:bCs/> before <w:b/> in the "Please note..." run properties.
$xml = $xml.Replace(
    '<w:r><w:rPr><w:b/></w:rPr><w:t>Please note',
    '<w:r><w:rPr><w:bCs/><w:b/></w:rPr><w:t>Please note')

# 2. Add w:val="0020" to the table's <w:tblLook/>.
$xml = $xml.Replace(
    '<w:tblLook w:firstRow="1" w:lastRow="0" w:firstColumn="0" w:lastColumn="0" w:noHBand="0" w:noVBand="0"/>',
    '<w:tblLook w:firstRow="1" w:lastRow="0" w:firstColumn="0" w:lastColumn="0" w:noHBand="0" w:noVBand="0" w:val="0020"/>')

# 3. Drop the explicit header-row <w:trPr><w:cnfStyle .../></w:trPr> and the
#    first header cell's <w:tcPr> (bottom border + bottom vAlign) -- these
#    are now supplied by the table style's new firstRow conditional format.
$xml = $xml.Replace(
    '<w:tr><w:trPr><w:cnfStyle w:firstRow="1" /></w:trPr><w:tc><w:tcPr><w:tcBorders><w:bottom w:val="single"/></w:tcBorders><w:vAlign w:val="bottom"/></w:tcPr>',
    '<w:tr><w:tc>')

# 4. Drop the same <w:tcPr> override from the other two header-row cells.
$xml = $xml.Replace(
    '</w:tc><w:tc><w:tcPr><w:tcBorders><w:bottom w:val="single"/></w:tcBorders><w:vAlign w:val="bottom"/></w:tcPr>',
    '</w:tc><w:tc>')

# 5. numbering.xml: the 9 levels of abstractNum 990 lose their explicit
#    <w:tabs><w:tab w:val="num" .../></w:tabs> and gain 240twips on the
#    hanging indent's w:left.
$xml = $xml.Replace(
    '<w:tabs><w:tab w:val="num" w:pos="0" /></w:tabs><w:ind w:left="480" w:hanging="480" />',
    '<w:ind w:left="720" w:hanging="480" />')
$xml = $xml.Replace(
    '<w:tabs><w:tab w:val="num" w:pos="720" /></w:tabs><w:ind w:left="1200" w:hanging="480" />',
    '<w:ind w:left="1440" w:hanging="480" />')
$xml = $xml.Replace(
    '<w:tabs><w:tab w:val="num" w:pos="1440" /></w:tabs><w:ind w:left="1920" w:hanging="480" />',
    '<w:ind w:left="2160" w:hanging="480" />')
$xml = $xml.Replace(
    '<w:tabs><w:tab w:val="num" w:pos="2160" /></w:tabs><w:ind w:left="2640" w:hanging="480" />',
    '<w:ind w:left="2880" w:hanging="480" />')
$xml = $xml.Replace(
    '<w:tabs><w:tab w:val="num" w:pos="2880" /></w:tabs><w:ind w:left="3360" w:hanging="480" />',
    '<w:ind w:left="3600" w:hanging="480" />')
$xml = $xml.Replace(
    '<w:tabs><w:tab w:val="num" w:pos="3600" /></w:tabs><w:ind w:left="4080" w:hanging="480" />',
    '<w:ind w:left="4320" w:hanging="480" />')
$xml = $xml.Replace(
    '<w:tabs><w:tab w:val="num" w:pos="4320" /></w:tabs><w:ind w:left="4800" w:hanging="480" />',
    '<w:ind w:left="5040" w:hanging="480" />')
$xml = $xml.Replace(
    '<w:tabs><w:tab w:val="num" w:pos="5040" /></w:tabs><w:ind w:left="5520" w:hanging="480" />',
    '<w:ind w:left="5760" w:hanging="480" />')
$xml = $xml.Replace(
    '<w:tabs><w:tab w:val="num" w:pos="5760" /></w:tabs><w:ind w:left="6240" w:hanging="480" />',
    '<w:ind w:left="6480" w:hanging="480" />')

# 6. styles.xml: give the "Table" style a firstRow conditional block that
#    reinstates the formatting removed from the table markup in step 3/4.
$oldTableStyleTail = @'
<w:tblCellMar>
        <w:top w:w="0" w:type="dxa" />
        <w:left w:w="108" w:type="dxa" />
        <w:bottom w:w="0" w:type="dxa" />
        <w:right w:w="108" w:type="dxa" />
      </w:tblCellMar>
    </w:tblPr></w:style>
'@
$newTableStyleTail = @'
<w:tblCellMar>
        <w:top w:w="0" w:type="dxa" />
        <w:left w:w="108" w:type="dxa" />
        <w:bottom w:w="0" w:type="dxa" />
        <w:right w:w="108" w:type="dxa" />
      </w:tblCellMar>
    </w:tblPr><w:tblStylePr w:type="firstRow"><w:tblPr><w:jc w:val="left"/><w:tblInd w:w="0" w:type="dxa"/></w:tblPr><w:trPr><w:jc w:val="left"/></w:trPr><w:tcPr><w:vAlign w:val="bottom"/><w:tcBorders><w:bottom w:val="single"/></w:tcBorders></w:tcPr></w:tblStylePr></w:style>
'@
$xml = $xml.Replace($oldTableStyleTail, $newTableStyleTail)

$d.WordOpenXML = $xml
Write-Output "applied"
